$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 16400.143
$ws.Range("I40").Value = 34767
$ws.Range("J40").Value = 2625
$ws.Range("K40").Value = 34767
$ws.Range("L40").Value = 2625
$ws.Range("M40").Value = -34592
$ws.Range("N40").Value = -2975
# Row 132
$ws.Range("H132").Value = 12901.23
$ws.Range("I132").Value = 1923.7858
$ws.Range("K132").Value = 5771.357400000001
$ws.Range("M132").Value = -3241.357400000001
# Row 137
$ws.Range("H137").Value = 2632.4934
$ws.Range("I137").Value = 911.8868
$ws.Range("J137").Value = 6777.591
$ws.Range("K137").Value = 2735.6604
$ws.Range("L137").Value = 20332.773
$ws.Range("M137").Value = -185.6603999999998
$ws.Range("N137").Value = -25432.773

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1896.9572
$ws.Range("I61").Value = 1074.122
$ws.Range("J61").Value = 3060.276
$ws.Range("K61").Value = 1074.122
$ws.Range("L61").Value = 3060.276
$ws.Range("M61").Value = -862.1220000000001
$ws.Range("N61").Value = -3484.276
# Row 74
$ws.Range("H74").Value = 2446.375
$ws.Range("I74").Value = 2196.0688
$ws.Range("J74").Value = 4866
$ws.Range("K74").Value = 2196.0688
$ws.Range("L74").Value = 4866
$ws.Range("M74").Value = -1322.0688
$ws.Range("N74").Value = -6614
# Row 77
$ws.Range("H77").Value = 2446.375
$ws.Range("I77").Value = 2196.0688
$ws.Range("J77").Value = 4866
$ws.Range("K77").Value = 10980.344
$ws.Range("L77").Value = 24330
$ws.Range("M77").Value = -6612.344000000001
$ws.Range("N77").Value = -33066
# Row 107
$ws.Range("H107").Value = 34590.4
$ws.Range("J107").Value = 34590.4
$ws.Range("L107").Value = 34590.4
$ws.Range("N107").Value = -42270.4
# Row 132
$ws.Range("H132").Value = 6580175.5
$ws.Range("I132").Value = 8197523
$ws.Range("K132").Value = 24592569
$ws.Range("M132").Value = -24590039
# Row 134
$ws.Range("H134").Value = 51800
$ws.Range("J134").Value = 51800
$ws.Range("L134").Value = 51800
$ws.Range("N134").Value = -61940
# Row 136
$ws.Range("H136").Value = 1896.9572
$ws.Range("I136").Value = 1074.122
$ws.Range("J136").Value = 3060.276
$ws.Range("K136").Value = 3222.366
$ws.Range("L136").Value = 9180.828
$ws.Range("M136").Value = -672.366
$ws.Range("N136").Value = -14280.828

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 166005.84
$ws.Range("I134").Value = 1329.069
$ws.Range("J134").Value = 249788.75
$ws.Range("K134").Value = 3987.207
$ws.Range("L134").Value = 749366.25
$ws.Range("M134").Value = -1452.207
$ws.Range("N134").Value = -754436.25

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2974.05
$ws.Range("I31").Value = 1290.7826
$ws.Range("J31").Value = 3476.8442
$ws.Range("K31").Value = 1290.7826
$ws.Range("L31").Value = 3476.8442
$ws.Range("M31").Value = -995.7826
$ws.Range("N31").Value = -4066.8442
# Row 34
$ws.Range("H34").Value = 2974.05
$ws.Range("I34").Value = 1290.7826
$ws.Range("J34").Value = 3476.8442
$ws.Range("K34").Value = 1290.7826
$ws.Range("L34").Value = 3476.8442
$ws.Range("M34").Value = -1088.7826
$ws.Range("N34").Value = -3880.8442
# Row 68
$ws.Range("H68").Value = 60000
$ws.Range("J68").Value = 60000
$ws.Range("L68").Value = 60000
$ws.Range("N68").Value = -61498
# Row 71
$ws.Range("H71").Value = 60000
$ws.Range("J71").Value = 60000
$ws.Range("L71").Value = 180000
$ws.Range("N71").Value = -187488
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
# Row 115
$ws.Range("H115").Value = 28234.25
$ws.Range("J115").Value = 28234.25
$ws.Range("L115").Value = 28234.25
$ws.Range("N115").Value = -30584.25
# Row 122
$ws.Range("H122").Value = 110799.09
$ws.Range("I122").Value = 172941.28
$ws.Range("J122").Value = 2050.25
$ws.Range("K122").Value = 518823.84
$ws.Range("L122").Value = 6150.75
$ws.Range("M122").Value = -516373.84
$ws.Range("N122").Value = -11050.75

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 6287.7144
$ws.Range("I3").Value = 2826.6667
$ws.Range("J3").Value = 8883.5
$ws.Range("K3").Value = 8480.000100000001
$ws.Range("L3").Value = 26650.5
$ws.Range("M3").Value = -8368.000100000001
$ws.Range("N3").Value = -26874.5
# Row 68
$ws.Range("H68").Value = 2878.1667
$ws.Range("I68").Value = 433.57144
$ws.Range("J68").Value = 6300.6
$ws.Range("K68").Value = 1300.71432
$ws.Range("L68").Value = 18901.8
$ws.Range("M68").Value = -489.71432
$ws.Range("N68").Value = -20523.8
# Row 70
$ws.Range("H70").Value = 4430.5
$ws.Range("I70").Value = 1573.1428
$ws.Range("J70").Value = 5969.077
$ws.Range("K70").Value = 4719.428400000001
$ws.Range("L70").Value = 17907.231
$ws.Range("M70").Value = -4404.428400000001
$ws.Range("N70").Value = -18537.231
# Row 71
$ws.Range("H71").Value = 2878.1667
$ws.Range("I71").Value = 433.57144
$ws.Range("J71").Value = 6300.6
$ws.Range("K71").Value = 3902.14296
$ws.Range("L71").Value = 56705.4
$ws.Range("M71").Value = 153.8570399999999
$ws.Range("N71").Value = -64817.4
# Row 73
$ws.Range("H73").Value = 4430.5
$ws.Range("I73").Value = 1573.1428
$ws.Range("J73").Value = 5969.077
$ws.Range("K73").Value = 4719.428400000001
$ws.Range("L73").Value = 17907.231
$ws.Range("M73").Value = -3627.428400000001
$ws.Range("N73").Value = -20091.231
# Row 113
$ws.Range("H113").Value = 3867.3547
$ws.Range("I113").Value = 8361.308000000001
$ws.Range("J113").Value = 621.7222
$ws.Range("K113").Value = 25083.924
$ws.Range("L113").Value = 1865.1666
$ws.Range("M113").Value = -22913.924
$ws.Range("N113").Value = -6205.1666
# Row 127
$ws.Range("H127").Value = 1020.75
$ws.Range("J127").Value = 1020.75
$ws.Range("L127").Value = 3062.25
$ws.Range("N127").Value = -12982.25
# Row 129
$ws.Range("H129").Value = 112550.41
$ws.Range("I129").Value = 188357.5
$ws.Range("J129").Value = 2285.5454
$ws.Range("K129").Value = 565072.5
$ws.Range("L129").Value = 6856.6362
$ws.Range("M129").Value = -560072.5
$ws.Range("N129").Value = -16856.6362
# Row 131
$ws.Range("H131").Value = 3715.55
$ws.Range("J131").Value = 1508.931
$ws.Range("L131").Value = 4526.793
$ws.Range("N131").Value = -14606.793
# Row 137
$ws.Range("H137").Value = 31258670
$ws.Range("I137").Value = 3117
$ws.Range("J137").Value = 45465740
$ws.Range("K137").Value = 9351
$ws.Range("L137").Value = 136397220
$ws.Range("M137").Value = -4251
$ws.Range("N137").Value = -136407420

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1789.7646
$ws.Range("I113").Value = 1820
$ws.Range("J113").Value = 1746.5714
$ws.Range("K113").Value = 1820
$ws.Range("L113").Value = 1746.5714
$ws.Range("M113").Value = 350
$ws.Range("N113").Value = -6086.5714

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2200.2222
$ws.Range("I7").Value = 1930.5
$ws.Range("J7").Value = 2739.6667
$ws.Range("K7").Value = 1930.5
$ws.Range("L7").Value = 2739.6667
$ws.Range("M7").Value = -1818.5
$ws.Range("N7").Value = -2963.6667
# Row 74
$ws.Range("H74").Value = 20707
$ws.Range("I74").Value = 11197
$ws.Range("J74").Value = 30217
$ws.Range("K74").Value = 11197
$ws.Range("L74").Value = 30217
$ws.Range("M74").Value = -10199
$ws.Range("N74").Value = -32213
# Row 77
$ws.Range("H77").Value = 20707
$ws.Range("I77").Value = 11197
$ws.Range("J77").Value = 30217
$ws.Range("K77").Value = 33591
$ws.Range("L77").Value = 90651
$ws.Range("M77").Value = -28599
$ws.Range("N77").Value = -100635
# Row 126
$ws.Range("H126").Value = 2200.2222
$ws.Range("I126").Value = 1930.5
$ws.Range("J126").Value = 2739.6667
$ws.Range("K126").Value = 5791.5
$ws.Range("L126").Value = 8219.000100000001
$ws.Range("M126").Value = -3321.5
$ws.Range("N126").Value = -13159.0001
# Row 132
$ws.Range("H132").Value = 2082.4648
$ws.Range("I132").Value = 1298.5962
$ws.Range("J132").Value = 4227.7896
$ws.Range("K132").Value = 3895.7886
$ws.Range("L132").Value = 12683.3688
$ws.Range("M132").Value = -1365.7886
$ws.Range("N132").Value = -17743.3688

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3246.6667
$ws.Range("I62").Value = 2700
$ws.Range("J62").Value = 3285.7144
$ws.Range("K62").Value = 2700
$ws.Range("L62").Value = 3285.7144
$ws.Range("M62").Value = -2076
$ws.Range("N62").Value = -4533.7144
# Row 65
$ws.Range("H65").Value = 3246.6667
$ws.Range("I65").Value = 2700
$ws.Range("J65").Value = 3285.7144
$ws.Range("K65").Value = 13500
$ws.Range("L65").Value = 16428.572
$ws.Range("M65").Value = -10380
$ws.Range("N65").Value = -22668.572
# Row 119
$ws.Range("H119").Value = 250034670
$ws.Range("J119").Value = 250034670
$ws.Range("L119").Value = 250034670
$ws.Range("N119").Value = -250044346
# Row 120
$ws.Range("H120").Value = 45416
$ws.Range("J120").Value = 45416
$ws.Range("L120").Value = 45416
$ws.Range("N120").Value = -55092
# Row 124
$ws.Range("H124").Value = 30183.857
$ws.Range("J124").Value = 30183.857
$ws.Range("L124").Value = 30183.857
$ws.Range("N124").Value = -40003.857
# Row 126
$ws.Range("H126").Value = 1280583.2
$ws.Range("I126").Value = 1472120.5
$ws.Range("K126").Value = 4416361.5
$ws.Range("M126").Value = -4413891.5
# Row 132
$ws.Range("H132").Value = 964.88525
$ws.Range("I132").Value = 450.07318
$ws.Range("J132").Value = 2020.25
$ws.Range("K132").Value = 1350.21954
$ws.Range("L132").Value = 6060.75
$ws.Range("M132").Value = 1179.78046
$ws.Range("N132").Value = -11120.75
# Row 136
$ws.Range("H136").Value = 14265.986
$ws.Range("I136").Value = 17414.88
$ws.Range("J136").Value = 1880.3334
$ws.Range("K136").Value = 52244.64
$ws.Range("L136").Value = 5641.0002
$ws.Range("M136").Value = -49694.64
$ws.Range("N136").Value = -10741.0002
